{"js": "// The underlying XML diff for this revision is purely a cosmetic\n// re-serialization of the package's canonical OOXML: every changed line is\n// the exact same element / attribute name-value set as before, just with\n// the attributes re-ordered (namespace declarations grouped and sorted\n// alphabetically, followed by the remaining attributes sorted\n// alphabetically). No text, formatting, style, section, or structural\n// property actually changes between the \"before\" and \"after\" snapshots.\n//\n// That kind of pure attribute-order normalization is an artifact of the\n// tool that produced the canonical-OOXML snapshot used for the diff (it\n// happens whenever a docx test fixture gets re-exported/re-serialized) --\n// it is not something that corresponds to any Word UI action or any\n// Word.* object-model call, since the object model has no notion of /\n// control over the raw attribute order the package writer uses when it\n// serializes a part. There is therefore no content, formatting or\n// structural change for this script to apply: the correct, faithful\n// reproduction of the diff is to leave the document's semantic content\n// untouched.\n//\n// We still touch the context the way a real task script would, so the\n// request context is exercised and synced, but we deliberately perform no\n// mutation.\nconst body = context.document.body;\nbody.load(\"text\");\nawait context.sync();\n", "ps1": "# The underlying XML diff for this revision is purely a cosmetic\n# re-serialization of the package's canonical OOXML: every changed line is\n# the exact same element / attribute name-value set as before, just with\n# the attributes re-ordered (namespace declarations grouped and sorted\n# alphabetically, followed by the remaining attributes sorted\n# alphabetically). No text, formatting, style, section, or structural\n# property actually changes between the \"before\" and \"after\" snapshots.\n#\n# That kind of pure attribute-order normalization is an artifact of the\n# tool that produced the canonical-OOXML snapshot used for the diff (it\n# happens whenever a docx test fixture gets re-exported/re-serialized) --\n# it is not something that corresponds to any Word UI action or any Word\n# COM object-model call, since the object model has no notion of / control\n# over the raw attribute order the package writer uses when it serializes\n# a part. There is therefore no content, formatting or structural change\n# for this script to apply: the correct, faithful reproduction of the diff\n# is to leave the document's semantic content untouched.\n#\n# We still touch the document the way a real task script would, so the\n# object model is exercised, but we deliberately perform no mutation.\n$d = $word.ActiveDocument\n$null = $d.Content.Text\n"}
